$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.953031599521637
$ws.Range("B1").Value = 1.484264731407166
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.367990374565125
